$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.671.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.417.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "533.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.70%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.413.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.97%  "

# Row 9
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.635"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.135"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.948.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.25%  "

# Row 16
$ws.Range("E16").Value = "  -2.61%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.402.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.333.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.55%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.986"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.74%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.90%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.58%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.43%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.92%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "688.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.53%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -16.94%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.95%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.107"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.91%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -12.20%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.388"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "28.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +29.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.903.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.81%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0403"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.37%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0634"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -16.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.127"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.73%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "

# Row 51
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.69%  "
